$wb = $excel.ActiveWorkbook

# --- survey sheet: add "hideInContents" column (O) ---
$survey = $wb.Worksheets.Item("survey")
$survey.Activate()
$survey.Range("O1").Value = "hideInContents"
$survey.Range("O2").Value = $true
$survey.Range("O10").Value = $true
$survey.Columns.Item(15).ColumnWidth = 13.5

# Update the selection on the survey sheet (per diff: C6 -> A10)
$survey.Range("A10").Select()

# --- settings sheet becomes the active tab ---
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
